# Adjust forecast parameters to increase weight of past data in smoothing
# (NH vaccinations have step variations that are probably one-time).

$wb = $excel.ActiveWorkbook

$sites = $wb.Worksheets.Item("SITES")
$source = $wb.Worksheets.Item("SOURCE")

# --- SITES sheet: update Ni (capacity) values -----------------------------
# Connecticut
$sites.Range("D2").Value = 131   # Medical facility
$sites.Range("D3").Value = 153   # Pharmacy
$sites.Range("D6").Value = 6     # Type 3

# Massachusetts
$sites.Range("D10").Value = 170  # Pharmacy
$sites.Range("D11").Value = 3    # Type 1

# Maine
$sites.Range("D16").Value = 70   # Medical facility
$sites.Range("D17").Value = 122  # Pharmacy
$sites.Range("D20").Value = 5    # Type 3
$sites.Range("D21").Value = 1    # Type 4

# Rhode Island
$sites.Range("D31").Value = 55   # Pharmacy
$sites.Range("D33").Value = 2    # Type 2
$sites.Range("D35").Value = 2    # Type 4 (was formula =ROUND(D29/2,0) -> now a plain value)

# Vermont
$sites.Range("D38").Value = 66   # Pharmacy

# --- SOURCE sheet: swap KV / source-document note -------------------------
$source.Range("A2").Value = "KV"
$source.Range("F2").Value = "See 2021-03-30 Vax Sites per State.docx"

# --- View / selection tidy-up ----------------------------------------------
$sites.Activate() | Out-Null
$sites.Range("D37").Select() | Out-Null

$source.Activate() | Out-Null
$source.Range("A3").Select() | Out-Null

$sites.Activate() | Out-Null
